# Update the Filtered_APAP_genes worksheet:
#  - Replace the ICAM1 gene row (row 12) with a new BAX gene row
#  - Replace the MIR122 gene row (row 16) with a new GPT gene row
#  - Update the sheet view (zoom level + active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filtered_APAP_genes")

# Row 12: was ICAM1 / 3383 / "...ICAM1 mRNA" -> now BAX / 959 / "...BAX mRNA"
$ws.Cells.Item(12, 4).Value = "BAX"
$ws.Cells.Item(12, 5).Value = 959
$ws.Cells.Item(12, 6).Value = "Acetaminophen results in increased expression of BAX mRNA"

# Row 16: was MIR122 / 406906 / "...MIR122 mRNA" -> now GPT / 4552 / "...GPT mRNA"
$ws.Cells.Item(16, 4).Value = "GPT"
$ws.Cells.Item(16, 5).Value = 4552
$ws.Cells.Item(16, 6).Value = "Acetaminophen results in increased expression of GPT mRNA"

# Update the view: zoom to 140% and select F16
$ws.Activate()
$excel.ActiveWindow.Zoom = 140
$ws.Range("F16").Select() | Out-Null
